# Fruta / hortaliza, semanal
# New weekly data point arrived for "Terminal La Palmera de La Serena - Membrillo":
# insert it as the new row 9, pushing the existing rows 9-18 down to 10-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 9; this shifts the old
# rows 9..18 down to 10..19 and keeps all their values intact.
$ws.Rows("9:9").Insert()

# Populate the freshly inserted row 9 with the new record.
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Terminal La Palmera de La Serena"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 44648
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100104
$ws.Range("H9").Value = "Frutos de pepita"
$ws.Range("I9").Value = 100104003
$ws.Range("J9").Value = "Membrillo"
$ws.Range("K9").Value = "Champion"
$ws.Range("L9").Value = "Especial"
$ws.Range("M9").Value = 16
$ws.Range("N9").Value = 335000
$ws.Range("O9").Value = 340000
$ws.Range("P9").Value = 337500
$ws.Range("Q9").Value = "$/bins (450 kilos)"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 750
$ws.Range("T9").Value = 450
